$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timing issue fix - keywords, updated tc1,2 in ubc01
# Replace the CasesTab (TC1/TC2) Neo4j query in B2 with the corrected
# version: the MATCH clauses were reordered/retimed and the trailing
# `Cohort` column was dropped from the RETURN clause.
$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n`tWHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in [ 'T3N0M1', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

# The shorter query text re-wraps to a shorter row, matching the other
# (already-updated) query rows.
$ws.Rows(2).RowHeight = 290

# Reset the view: land on A2 with B2 selected (was parked on C4 before).
$ws.Range("A2").Select()
$ws.Range("B2").Select()
